# Sync preset values with updated Presets.xlsx
# Applies the numeric preset changes described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Presets")

# --- Intensity Preset (row 3 multipliers) ---
$ws.Range("C3").Value = 1.5    # Subtle
$ws.Range("E3").Value = 0.8    # Dramatic
$ws.Range("F3").Value = 0.5    # Cinematic
$ws.Range("G3").Value = 0.3    # Epic

# Last Stand intensity base (timeScale)
$ws.Range("B11").Value = 0.3

# --- Duration Preset (row 16 multipliers) ---
$ws.Range("C16").Value = 0.35  # Very Short
$ws.Range("F16").Value = 1.35  # Long
$ws.Range("G16").Value = 1.7   # Extended

# Last Stand duration base
$ws.Range("B24").Value = 4.0

# --- Cooldown Preset ---
$ws.Range("B36").Value = 30    # Last Enemy cooldown base
$ws.Range("B37").Value = 90    # Last Stand cooldown base

# --- Chance Preset ---
$ws.Range("B46").Value = 0.3   # Dismemberment chance base
